# A new weekly observation was added to the price table.
# This inserts a new row at position 35 (pushing the existing rows
# 35..152 down to 36..153) and fills it with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(35).Insert()

$ws.Range("A35").Value = 4
$ws.Range("B35").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C35").Value = "Los Lagos"
$ws.Range("D35").Value = 44624
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = 100112009
$ws.Range("G35").Value = "Acelga"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 80
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 10000
$ws.Range("N35").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O35").Value = "Región de La Araucanía"
$ws.Range("P35").Value = 833
$ws.Range("Q35").Value = 12
$ws.Range("R35").Value = "Hortaliza"
